$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new data rows above the existing (last) employee row (row 16),
# shifting it down to row 19, along with the signature rows below.
$ws.Rows("16:18").Insert()

# Seed the three new rows with the formatting of the row that is now last
# (row 19 still carries the original "closing" table style at this point).
$ws.Range("B19:J19").Copy($ws.Range("B16:J16"))
$ws.Range("B19:J19").Copy($ws.Range("B17:J17"))
$ws.Range("B19:J19").Copy($ws.Range("B18:J18"))

# Fill in the three new employee rows.
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "9096736"
$ws.Range("D16").Value = "JAIME WALTER RODRIGUEZ ALVAREZ"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 60000
$ws.Range("G16").Value = 1500000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73100724"
$ws.Range("D17").Value = "FERNANDO CABALLERO DIAZ GRANADOS"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 160000
$ws.Range("G17").Value = 4000000

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "23140579"
$ws.Range("D18").Value = "ANA LEONOR ROMERO MOLINA"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 60000
$ws.Range("G18").Value = 1500000

# Update the summary figures at the top of the sheet.
$ws.Range("E11").Value = 306000
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 2

# Column D needs to be a bit wider to fit the longest new name.
$ws.Columns("D:D").ColumnWidth = 37.7265625

Write-Output "edit applied"
